$d = $word.ActiveDocument

$replacements = @(
    @{old = "13×32=416"; new = "89×18=1602"},
    @{old = "75×72=5400"; new = "33×89=2937"},
    @{old = "69×38=2622"; new = "41×65=2665"},
    @{old = "97×76=7372"; new = "50×80=4000"},
    @{old = "74×79=5846"; new = "75×49=3675"},
    @{old = "72×36=2592"; new = "99×80=7920"},
    @{old = "82×88=7216"; new = "57×16=912"},
    @{old = "28×54=1512"; new = "88×20=1760"},
    @{old = "82×49=4018"; new = "41×25=1025"},
    @{old = "82×32=2624"; new = "16×61=976"},
    @{old = "12×19=228"; new = "80×87=6960"},
    @{old = "17×48=816"; new = "29×45=1305"},
    @{old = "82×35=2870"; new = "83×48=3984"},
    @{old = "77×15=1155"; new = "20×85=1700"},
    @{old = "56×19=1064"; new = "76×92=6992"},
    @{old = "56×59=3304"; new = "22×45=990"},
    @{old = "63×38=2394"; new = "62×69=4278"},
    @{old = "17×98=1666"; new = "18×25=450"},
    @{old = "33×15=495"; new = "73×46=3358"},
    @{old = "21×38=798"; new = "35×65=2275"},
    @{old = "25×51=1275"; new = "18×52=936"},
    @{old = "78×60=4680"; new = "70×24=1680"},
    @{old = "31×30=930"; new = "32×94=3008"},
    @{old = "56×84=4704"; new = "42×69=2898"},
    @{old = "25×80=2000"; new = "63×58=3654"}
)

foreach ($r in $replacements) {
    $find = $d.Content
    $find.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
